# TalDoor_BOM.xlsx update
# The C15 capacitor (row 11 of the "TalDoor" BOM sheet) is being swapped for
# a different part: value 1000uF -> 470uF, with a new Digikey part number
# and updated unit price (which ripples into the "Price" extended-cost
# formula and the running total at G36 automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TalDoor")

# Value column (C) for the C15, row -- was "1000uF"
$ws.Range("C11").Value = "470uF"

# Digikey Part # column (F) for the C15 row -- was "1572-1665-ND"
$ws.Range("F11").Value = "493-11709-1-ND"

# Unit price column (G) for the C15 row -- was 0.283
$ws.Range("G11").Value = 0.311

# Best-effort: nudge the saved window scroll position (cosmetic, yWindow
# 6000 -> 6600 in bookViews) via the Window object, if the host surfaces it.
try {
    $win = $excel.ActiveWindow
    $win.Top = 6600
} catch {
}
